$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# Row 169 is an existing skill row ("PointSelf"/"OnAdd"/"GetDescript" style,
# single-target buff-on-add skill) whose cell formatting (styles) is an exact
# match for the new skill row we need to insert at row 184. Copy it and
# insert a fresh row above the old row 184 so every row from 184..200 shifts
# down to 185..201, carrying the correct per-cell styles along with it.
$ws.Rows.Item(169).Copy()
$ws.Rows.Item(184).Insert()

# New skill: 55900052 "陪练" (sparring partner) - on add, taunts the closest
# friendly unit within range.
$ws.Rows.Item(184).RowHeight = 120

$ws.Cells.Item(184, 1).Value = 55900052
$ws.Cells.Item(184, 2).Value = "陪练"
$ws.Cells.Item(184, 3).Value = "特殊"
$ws.Cells.Item(184, 4).Value = "NFR"
$ws.Cells.Item(184, 5).Value = 15
$ws.Cells.Item(184, 6).Value = "true"
$ws.Cells.Item(184, 7).Value = ""
$ws.Cells.Item(184, 8).Value = "foreach(IMonster mon in s.Map.GetRangeMonster(s.IsLeft,sp.Target,sp.Shape,sp.Range,s.Position).FilterId(s.Id).SortDistance(true).Top(1)) mon.AddSkill(55100008,100);"
$ws.Cells.Item(184, 9).Value = ""
$ws.Cells.Item(184, 10).Value = ""
$ws.Cells.Item(184, 11).Value = ""
$ws.Cells.Item(184, 12).Value = ""
$ws.Cells.Item(184, 13).Value = ""
$ws.Cells.Item(184, 14).Value = ""
$ws.Cells.Item(184, 15).Value = ""
$ws.Cells.Item(184, 17).Value = "true"
$ws.Cells.Item(184, 18).Value = "true"
$ws.Cells.Item(184, 19).Value = "给予范围内最近友方单位嘲讽"
$ws.Cells.Item(184, 20).Value = ""
$ws.Cells.Item(184, 21).Value = ""
$ws.Cells.Item(184, 22).Value = ""
$ws.Cells.Item(184, 23).Value = ""
$ws.Cells.Item(184, 24).Value = 5
$ws.Cells.Item(184, 25).Value = "peilian"
$ws.Cells.Item(184, 26).Value = ""
